$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "isnegative"
$ws.Range("G1").Value = "shift"
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0

$ws.Range("F2").Select()
